$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text (e.g. "613.71") that
# Excel would otherwise auto-convert to a real number on assignment. Force
# the whole data range to Text format first so values are written as
# strings (matching the inline-string cells in the original workbook),
# then restore the default "Normal" style so no stray formatting remains.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "70.568.88"
$ws.Cells.Item(2, 5).Value = "  +2.38%  "
$ws.Cells.Item(3, 4).Value = "3.967.02"
$ws.Cells.Item(3, 5).Value = "  +3.15%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).Value = "613.71"
$ws.Cells.Item(5, 5).Value = "  +1.76%  "
$ws.Cells.Item(6, 4).Value = "171.41"
$ws.Cells.Item(6, 5).Value = "  +4.99%  "
$ws.Cells.Item(7, 4).Value = "3.964.25"
$ws.Cells.Item(7, 5).Value = "  +3.16%  "
$ws.Cells.Item(8, 5).Value = "  -0.04%  "
$ws.Cells.Item(9, 4).Value = "0.540"
$ws.Cells.Item(9, 5).Value = "  +1.55%  "
$ws.Cells.Item(10, 5).Value = "  +2.86%  "
$ws.Cells.Item(11, 4).Value = "6.51"
$ws.Cells.Item(11, 5).Value = "  +3.00%  "
$ws.Cells.Item(12, 4).Value = "0.474"
$ws.Cells.Item(12, 5).Value = "  +3.08%  "
$ws.Cells.Item(13, 4).Value = "0.0000260"
$ws.Cells.Item(13, 5).Value = "  +6.38%  "
$ws.Cells.Item(14, 4).Value = "38.80"
$ws.Cells.Item(14, 5).Value = "  +5.09%  "
$ws.Cells.Item(15, 4).Value = "4.626.64"
$ws.Cells.Item(15, 5).Value = "  +3.20%  "
$ws.Cells.Item(16, 4).Value = "3.941.84"
$ws.Cells.Item(16, 5).Value = "  +2.57%  "
$ws.Cells.Item(17, 4).Value = "70.383.34"
$ws.Cells.Item(17, 5).Value = "  +1.87%  "
$ws.Cells.Item(18, 4).Value = "7.74"
$ws.Cells.Item(18, 5).Value = "  +1.67%  "
$ws.Cells.Item(19, 4).Value = "18.21"
$ws.Cells.Item(19, 5).Value = "  +6.08%  "
$ws.Cells.Item(20, 5).Value = "  -1.35%  "
$ws.Cells.Item(21, 4).Value = "11.04"
$ws.Cells.Item(21, 5).Value = "  -3.72%  "
$ws.Cells.Item(22, 4).Value = "501.68"
$ws.Cells.Item(22, 5).Value = "  +2.97%  "
$ws.Cells.Item(23, 4).Value = "0.747"
$ws.Cells.Item(23, 5).Value = "  +3.70%  "
$ws.Cells.Item(24, 4).Value = "0.0000169"
$ws.Cells.Item(24, 5).Value = "  +7.07%  "
$ws.Cells.Item(25, 4).Value = "86.09"
$ws.Cells.Item(25, 5).Value = "  +2.18%  "
$ws.Cells.Item(26, 5).Value = "  +2.17%  "
$ws.Cells.Item(27, 4).Value = "12.54"
$ws.Cells.Item(27, 5).Value = "  +3.00%  "
$ws.Cells.Item(28, 4).Value = "10.30"
$ws.Cells.Item(28, 5).Value = "  +2.78%  "
$ws.Cells.Item(29, 5).Value = "  +0.26%  "
$ws.Cells.Item(30, 4).Value = "3.03"
$ws.Cells.Item(30, 5).Value = "  +1.87%  "
$ws.Cells.Item(31, 4).Value = "4.117.95"
$ws.Cells.Item(31, 5).Value = "  +2.86%  "
$ws.Cells.Item(32, 5).Value = "  +1.77%  "
$ws.Cells.Item(33, 4).Value = "7.94"
$ws.Cells.Item(33, 5).Value = "  +0.00%  "
$ws.Cells.Item(34, 4).Value = "32.56"
$ws.Cells.Item(34, 5).Value = "  +0.78%  "
$ws.Cells.Item(35, 4).Value = "3.932.58"
$ws.Cells.Item(35, 5).Value = "  +3.76%  "
$ws.Cells.Item(36, 5).Value = "  +1.61%  "
$ws.Cells.Item(37, 4).Value = "6.22"
$ws.Cells.Item(37, 5).Value = "  +5.78%  "
$ws.Cells.Item(38, 5).Value = "  +1.80%  "
$ws.Cells.Item(39, 4).Value = "0.142"
$ws.Cells.Item(39, 5).Value = "  +1.12%  "
$ws.Cells.Item(40, 4).Value = "3.29"
$ws.Cells.Item(40, 5).Value = "  +10.33%  "
$ws.Cells.Item(41, 4).Value = "0.999"
$ws.Cells.Item(41, 5).Value = "  -0.02%  "
$ws.Cells.Item(42, 4).Value = "0.329"
$ws.Cells.Item(42, 5).Value = "  +3.07%  "
$ws.Cells.Item(43, 4).Value = "2.10"
$ws.Cells.Item(43, 5).Value = "  +5.97%  "
$ws.Cells.Item(44, 4).Value = "444.29"
$ws.Cells.Item(45, 4).Value = "48.40"
$ws.Cells.Item(45, 5).Value = "  -0.25%  "
$ws.Cells.Item(46, 4).Value = "8.68"
$ws.Cells.Item(46, 5).Value = "  +3.35%  "
$ws.Cells.Item(47, 5).Value = "  +0.01%  "
$ws.Cells.Item(48, 4).Value = "0.000279"
$ws.Cells.Item(48, 5).Value = "  +23.92%  "
$ws.Cells.Item(49, 2).Value = "Arweave"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(49, 4).Value = "41.12"
$ws.Cells.Item(49, 5).Value = "  +4.84%  "
$ws.Cells.Item(50, 2).Value = "VeChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(50, 4).Value = "0.0370"
$ws.Cells.Item(50, 5).Value = "  +3.41%  "
$ws.Cells.Item(51, 4).Value = "143.51"
$ws.Cells.Item(51, 5).Value = "  +0.63%  "

$priceRange.Style = "Normal"
